$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update G2:G5 ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("G2").Value = 49
$wsProd.Range("G3").Value = 21
$wsProd.Range("G4").Value = 35
$wsProd.Range("G5").Value = 70

# Column H holds empty-string placeholder cells; re-assert blank values so
# the sheet round-trips through COM without picking up a stray shared-string
# index for these previously-empty cells.
$wsProd.Range("H2:H11").Value = ""

# --- ForecastedAverageDemand sheet: update B9:E11 ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvg.Range("B9").Value = 70
$wsAvg.Range("C9").Value = 30
$wsAvg.Range("D9").Value = 50
$wsAvg.Range("E9").Value = 100

$wsAvg.Range("B10").Value = 70
$wsAvg.Range("C10").Value = 30
$wsAvg.Range("D10").Value = 50
$wsAvg.Range("E10").Value = 100

$wsAvg.Range("B11").Value = 70
$wsAvg.Range("C11").Value = 30
$wsAvg.Range("D11").Value = 50
$wsAvg.Range("E11").Value = 100

# --- ForcastedStandardDeviation sheet: update B9:E11 ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Range("B9").Value = 7.166424999999998
$wsStd.Range("C9").Value = 3.071324999999999
$wsStd.Range("D9").Value = 5.118874999999999
$wsStd.Range("E9").Value = 10.23775

$wsStd.Range("B10").Value = 8.1997825
$wsStd.Range("C10").Value = 3.5141925
$wsStd.Range("D10").Value = 5.856987499999999
$wsStd.Range("E10").Value = 11.713975

$wsStd.Range("B11").Value = 9.129804249999998
$wsStd.Range("C11").Value = 3.912773249999999
$wsStd.Range("D11").Value = 6.521288749999998
$wsStd.Range("E11").Value = 13.0425775
